$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 147, shifting existing rows 147-168 down to 149-170.
$ws.Rows("147:148").Insert()

# Fill the two newly inserted rows with this weeks price data.
# Row 147
$ws.Cells.Item(147,1).Value = 3
$ws.Cells.Item(147,2).Value = 'Femacal de La Calera'
$ws.Cells.Item(147,3).Value = 'Coquimbo'
$ws.Cells.Item(147,4).Value = 44449
$ws.Cells.Item(147,5).Value = 5
$ws.Cells.Item(147,6).Value = 100112013
$ws.Cells.Item(147,7).Value = 'Alcachofa'
$ws.Cells.Item(147,8).Value = 'Argentina(o)'
$ws.Cells.Item(147,9).Value = 'Primera'
$ws.Cells.Item(147,10).Value = 135
$ws.Cells.Item(147,11).Value = 10000
$ws.Cells.Item(147,12).Value = 11000
$ws.Cells.Item(147,13).Value = 10444
$ws.Cells.Item(147,14).Value = '$/caja 50 unidades'
$ws.Cells.Item(147,15).Value = 'Provincia de Limarí'
$ws.Cells.Item(147,16).Value = 209
$ws.Cells.Item(147,17).Value = 50
$ws.Cells.Item(147,18).Value = 'Hortaliza'

# Row 148
$ws.Cells.Item(148,1).Value = 3
$ws.Cells.Item(148,2).Value = 'Femacal de La Calera'
$ws.Cells.Item(148,3).Value = 'Coquimbo'
$ws.Cells.Item(148,4).Value = 44449
$ws.Cells.Item(148,5).Value = 5
$ws.Cells.Item(148,6).Value = 100112013
$ws.Cells.Item(148,7).Value = 'Alcachofa'
$ws.Cells.Item(148,8).Value = 'Española'
$ws.Cells.Item(148,9).Value = 'Extra'
$ws.Cells.Item(148,10).Value = 125
$ws.Cells.Item(148,11).Value = 11500
$ws.Cells.Item(148,12).Value = 12000
$ws.Cells.Item(148,13).Value = 11760
$ws.Cells.Item(148,14).Value = '$/caja 30 unidades'
$ws.Cells.Item(148,15).Value = 'Provincia de Limarí'
$ws.Cells.Item(148,16).Value = 392
$ws.Cells.Item(148,17).Value = 30
$ws.Cells.Item(148,18).Value = 'Hortaliza'

